$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

$ws.Cells.Item($row, 1).Value  = "Philippines"
$ws.Cells.Item($row, 2).Value  = "philippines"
$ws.Cells.Item($row, 3).Value  = "Panay River Basin"
$ws.Cells.Item($row, 4).Value  = "panay"
$ws.Cells.Item($row, 5).Value  = "Dao Bridge"
$ws.Cells.Item($row, 6).Value  = "G5369"
$ws.Cells.Item($row, 7).Value  = "primary"

# Column H holds a literal date-like string (e.g. "2025-10-29"), not an
# Excel date serial. Force the cell to Text before assignment so Excel's
# auto-detection doesn't convert it to a date, then restore the cell's
# style to "Normal" so no extra number-format styling is left behind.
$ws.Cells.Item($row, 8).NumberFormat = "@"
$ws.Cells.Item($row, 8).Value = "2025-10-29"
$ws.Cells.Item($row, 8).Style = "Normal"

$ws.Cells.Item($row, 9).Value  = 3
$ws.Cells.Item($row, 10).Value = 11.4249999999999
$ws.Cells.Item($row, 11).Value = 122.7249999999997
$ws.Cells.Item($row, 12).Value = 5
$ws.Cells.Item($row, 13).Value = 864.7614412809821
$ws.Cells.Item($row, 14).Value = "LOW"
$ws.Cells.Item($row, 15).Value = 603.6038567117938
$ws.Cells.Item($row, 16).Value = 864.7614412809821
$ws.Cells.Item($row, 17).Value = 50
$ws.Cells.Item($row, 18).Value = 0
$ws.Cells.Item($row, 19).Value = 0
$ws.Cells.Item($row, 20).Value = 243.98828125
$ws.Cells.Item($row, 21).Value = 252.8243713378906
$ws.Cells.Item($row, 22).Value = 151.3125
$ws.Cells.Item($row, 23).Value = 445.3984375
$ws.Cells.Item($row, 24).Value = 212.97265625
$ws.Cells.Item($row, 25).Value = 270.123046875
$ws.Cells.Item($row, 26).Value = $false
$ws.Cells.Item($row, 27).Value = -71.78548098900235
